$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 — duplicate of the existing Sharjah / Sunrisers match (currently row 3)
$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " November 03 2020"
$ws.Range("C4").Value = "Sunrisers won by 10 wickets (with 17 balls remaining)"
$ws.Range("D4").Value = "Mumbai Indians"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "Nathan Coulter-Nile "

# totalRuns/totalBalls/total4s/total6s/sr are stored as text (numberStoredAsText) —
# force a text format before assigning so the numeric-looking strings aren't
# silently coerced into real numbers.
$ws.Range("G4:K4").NumberFormat = "@"
$ws.Range("G4").Value = "1"
$ws.Range("H4").Value = "3"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "33.33"

# New row 5 — duplicate of the existing Dubai (DSC) / Kings XI match (currently row 2)
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 18 2020"
$ws.Range("C5").Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Range("D5").Value = "Mumbai Indians"
$ws.Range("E5").Value = "Kings XI Punjab"
$ws.Range("F5").Value = "Nathan Coulter-Nile "

$ws.Range("G5:K5").NumberFormat = "@"
$ws.Range("G5").Value = "24"
$ws.Range("H5").Value = "12"
$ws.Range("I5").Value = "4"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "200.00"
